# "Updated All Menu & Help"
# The product catalog's second row (key "AmaSearch") listed "Guitar" as its
# value; update it to "Shoes".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AmazonDataSheet")

$ws.Range("B2").Value = "Shoes"

# Reflect the zoomed-in view and final cell selection left behind after the
# edit (matches the sheetView state captured on save).
$win = $excel.ActiveWindow
$win.Zoom = 120
[void]$ws.Range("C11").Select()
